$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.733.82'
$ws.Range("E2").Value = '  +0.87%  '

$ws.Range("D3").Value = '1.657.62'
$ws.Range("E3").Value = '  +1.00%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.000'
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.85'
$ws.Range("E6").Value = '  -0.27%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3809'
$ws.Range("E7").Value = '  +0.27%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3619'
$ws.Range("E8").Value = '  -0.37%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.24'
$ws.Range("E9").Value = '  -0.72%  '

$ws.Range("E10").Value = '  +0.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.240'
$ws.Range("E11").Value = '  +0.22%  '

$ws.Range("E12").Value = '  -0.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.64'
$ws.Range("E13").Value = '  +0.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.523'
$ws.Range("E14").Value = '  +0.78%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.434'
$ws.Range("E15").Value = '  +0.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001233'
$ws.Range("E16").Value = '  -0.92%  '

$ws.Range("D17").Value = '1.642.55'
$ws.Range("E17").Value = '  +0.47%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '97.59'
$ws.Range("E18").Value = '  +2.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06990'
$ws.Range("E19").Value = '  +0.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.834'
$ws.Range("E20").Value = '  +3.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.68'
$ws.Range("E21").Value = '  +0.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.86'
$ws.Range("E23").Value = '  +2.22%  '

$ws.Range("D24").Value = '23.726.05'
$ws.Range("E24").Value = '  +0.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.512'
$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.045'
$ws.Range("E26").Value = '  -0.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.28'
$ws.Range("E27").Value = '  +0.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.99'
$ws.Range("E28").Value = '  +0.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.208'
$ws.Range("E29").Value = '  -1.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.56'
$ws.Range("E30").Value = '  +0.88%  '

$ws.Range("D31").Value = '1.836.59'
$ws.Range("E31").Value = '  +1.30%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.971'
$ws.Range("E32").Value = '  +4.94%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.188'
$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02817'
$ws.Range("E36").Value = '  +1.71%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2527'
$ws.Range("E37").Value = '  +1.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.115'
$ws.Range("E38").Value = '  +1.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.08791'
$ws.Range("E39").Value = '  +0.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.07056'
$ws.Range("E40").Value = '  -1.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.08'
$ws.Range("E41").Value = '  +7.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7040'
$ws.Range("E42").Value = '  -0.66%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.335'
$ws.Range("E43").Value = '  -0.69%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.05'
$ws.Range("E44").Value = '  +1.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6534'
$ws.Range("E45").Value = '  -0.43%  '

$ws.Range("E46").Value = '  +1.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  +0.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.978'
$ws.Range("E48").Value = '  +0.19%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07950'
$ws.Range("E49").Value = '  -0.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.15'
$ws.Range("E50").Value = '  +0.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.185'
$ws.Range("E51").Value = '  -0.95%  '

# Row 34/35: swap ImmutableX and FraxShare entries with updated data
$ws.Range("B34").Value = 'FraxShare'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.04'
$ws.Range("E34").Value = '  +5.07%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.064'
$ws.Range("E35").Value = '  +0.81%  '

